$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10840.4886431561
$ws.Range("F2").Value = 7.77819044789157

$ws.Range("C3").Value = 10199.9197492142
$ws.Range("F3").Value = 332.841584136756

$ws.Range("C4").Value = 7164.45359137731
$ws.Range("F4").Value = 179.382814163594

$ws.Range("C5").Value = 6981.74653789196
$ws.Range("F5").Value = 159.88138561046

$ws.Range("C6").Value = 10417.3390991584
$ws.Range("D6").Value = 11232.26
$ws.Range("F6").Value = 262.704367485383

$ws.Range("C7").Value = 10688.0522990171
$ws.Range("D7").Value = 11232.26
$ws.Range("F7").Value = 289.347973695107

$ws.Range("D8").Value = 11232.26
$ws.Range("F8").Value = 277.399569870025

$ws.Range("D9").Value = 11232.26
$ws.Range("F9").Value = 274.852408290793

$ws.Range("D10").Value = 11232.26
$ws.Range("F10").Value = 252.20020961704

$ws.Range("D11").Value = 11232.26
$ws.Range("F11").Value = 111.255279278

$ws.Range("D12").Value = 11232.26
$ws.Range("F12").Value = 115.199757931993

$ws.Range("D13").Value = 11232.26
$ws.Range("F13").Value = 151.780445675196

$ws.Range("D14").Value = 11232.26
$ws.Range("F14").Value = 343.256452553078

$ws.Range("D15").Value = 11232.26
$ws.Range("F15").Value = 350.750943286457
